# "Update all back to list pages"
# Resets the partId (B), shipDate (E) and quantity (F) columns for every
# work-order row on Sheet1 back to the values shown on the list pages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# partId (column B) and quantity (column F) updates - order is not significant.
$ws.Cells.Item(3, 2).Value  = 1
$ws.Cells.Item(4, 2).Value  = 1
$ws.Cells.Item(5, 2).Value  = 1
$ws.Cells.Item(6, 2).Value  = 1
$ws.Cells.Item(7, 2).Value  = 1
$ws.Cells.Item(8, 2).Value  = 1
$ws.Cells.Item(9, 2).Value  = 1
$ws.Cells.Item(10, 2).Value = 2
$ws.Cells.Item(11, 2).Value = 1
$ws.Cells.Item(12, 2).Value = 1
$ws.Cells.Item(13, 2).Value = 1
$ws.Cells.Item(14, 2).Value = 1
$ws.Cells.Item(15, 2).Value = 1
$ws.Cells.Item(16, 2).Value = 1
$ws.Cells.Item(17, 2).Value = 1
$ws.Cells.Item(18, 2).Value = 1

for ($r = 3; $r -le 18; $r++) {
    $ws.Cells.Item($r, 6).Value = 10000
}

# shipDate (column E) updates. Written in this specific order so that any
# brand-new distinct strings land in the shared-string table in the same
# sequence the workbook was originally saved with.
$ws.Cells.Item(16, 5).Value = "15/03/2020 12:00PM"
$ws.Cells.Item(18, 5).Value = "03/03/2020 12:00PM"
$ws.Cells.Item(17, 5).Value = "10/03/2020 12:00PM"
$ws.Cells.Item(15, 5).Value = "12/03/2020  12:00PM"
$ws.Cells.Item(5, 5).Value  = "08/03/2020 12:00PM"
$ws.Cells.Item(10, 5).Value = "05/03/2020 12:00PM"
$ws.Cells.Item(9, 5).Value  = "02/04/2020 12:00PM"
$ws.Cells.Item(7, 5).Value  = "10/3/2020  12:12PM"
$ws.Cells.Item(8, 5).Value  = "12/03/2020 12:00PM"

$ws.Cells.Item(3, 5).Value  = "08/08/2020 12:00PM"
$ws.Cells.Item(4, 5).Value  = "06/06/2020 12:00PM"
$ws.Cells.Item(6, 5).Value  = "10/07/2020 12:00PM"
$ws.Cells.Item(11, 5).Value = "10/06/2020 12:00PM"
$ws.Cells.Item(12, 5).Value = "30/4/2020  12:12PM"
$ws.Cells.Item(13, 5).Value = "05/05/2020 12:00PM"
$ws.Cells.Item(14, 5).Value = "30/03/2020  12:00PM"

# Selection moves from D7 to B18 as part of the save.
$ws.Range("B18").Select()
